$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to stay
# Text (matching the original inline-string/text cell contents) because
# Excel would otherwise auto-convert them into numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '56.634.47'
$ws.Range("E2").Value = '  -3.89%  '
$ws.Range("D3").Value = '2.361.97'
$ws.Range("E3").Value = '  -6.34%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '514.24'
$ws.Range("E5").Value = '  -4.12%  '
$ws.Range("D6").Value = '127.95'
$ws.Range("E6").Value = '  -5.70%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").Value = '0.554'
$ws.Range("E8").Value = '  -2.25%  '
$ws.Range("D9").Value = '2.377.67'
$ws.Range("E9").Value = '  -5.67%  '
$ws.Range("D10").Value = '0.0960'
$ws.Range("E10").Value = '  -3.76%  '
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("E12").Value = '  -8.44%  '
$ws.Range("D13").Value = '0.317'
$ws.Range("E13").Value = '  -5.60%  '
$ws.Range("D14").Value = '2.779.04'
$ws.Range("E14").Value = '  -6.45%  '
$ws.Range("D15").Value = '56.542.51'
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("D16").Value = '21.48'
$ws.Range("E16").Value = '  -4.71%  '
$ws.Range("E17").Value = '  -4.14%  '
$ws.Range("D18").Value = '2.341.25'
$ws.Range("E18").Value = '  -7.23%  '
$ws.Range("D19").Value = '10.29'
$ws.Range("E19").Value = '  -4.17%  '
$ws.Range("D20").Value = '310.03'
$ws.Range("E20").Value = '  -3.94%  '
$ws.Range("D21").Value = '4.03'
$ws.Range("E21").Value = '  -5.23%  '
$ws.Range("D22").Value = '6.11'
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '64.67'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '0.391'
$ws.Range("E26").Value = '  -4.48%  '
$ws.Range("D27").Value = '2.464.57'
$ws.Range("E27").Value = '  -6.79%  '
$ws.Range("E28").Value = '  -4.82%  '
$ws.Range("D29").Value = '7.17'
$ws.Range("E29").Value = '  -4.89%  '
$ws.Range("D30").Value = '173.26'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").Value = '1.67'
$ws.Range("E31").Value = '  -5.21%  '
$ws.Range("D32").Value = '0.0₃0717'
$ws.Range("E32").Value = '  -6.69%  '
$ws.Range("D33").Value = '6.10'
$ws.Range("E33").Value = '  -4.35%  '
$ws.Range("E34").Value = '  -7.18%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '17.59'
$ws.Range("E37").Value = '  -3.71%  '
$ws.Range("E38").Value = '  -6.22%  '
$ws.Range("D39").Value = '3.73'
$ws.Range("E39").Value = '  -7.14%  '
$ws.Range("D40").Value = '0.801'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("D41").Value = '35.49'
$ws.Range("E41").Value = '  -3.58%  '
$ws.Range("D42").Value = '1.43'
$ws.Range("E42").Value = '  -6.62%  '
$ws.Range("D43").Value = '3.33'
$ws.Range("E43").Value = '  -4.96%  '
$ws.Range("D44").Value = '4.89'
$ws.Range("E44").Value = '  -4.23%  '
$ws.Range("D45").Value = '122.72'
$ws.Range("E45").Value = '  -7.21%  '
$ws.Range("D46").Value = '0.570'
$ws.Range("E46").Value = '  -5.20%  '
$ws.Range("D47").Value = '252.49'
$ws.Range("E47").Value = '  -10.14%  '
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").Value = '0.0488'
$ws.Range("E49").Value = '  -4.52%  '
$ws.Range("D50").Value = '0.0208'
$ws.Range("E50").Value = '  -5.69%  '
$ws.Range("D51").Value = '16.68'
$ws.Range("E51").Value = '  -6.33%  '

# Restore the default (Normal) style on those cells so no stray
# number-format/style is left attached to them.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
